# "Fruta / hortaliza, semanal" — weekly refresh of the Repollo price sheet.
# A new week's record is inserted at row 68 (the data series is stored
# newest-first-ish, interleaved), pushing every existing record below it
# down by one row (old row 68 -> new row 69, ..., old row 116 -> new row 117).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 68, shifting rows 68:116 down to 69:117.
$ws.Rows(68).Insert()

# Populate the newly inserted row 68 with this week's record.
$ws.Range("A68").Value = 7
$ws.Range("B68").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C68").Value = "Ñuble"
$ws.Range("D68").Value = 44447
$ws.Range("E68").Value = 16
$ws.Range("F68").Value = 100112006
$ws.Range("G68").Value = "Repollo"
$ws.Range("H68").Value = "Crespo record"
$ws.Range("I68").Value = "Primera"
$ws.Range("J68").Value = 300
$ws.Range("K68").Value = 700
$ws.Range("L68").Value = 750
$ws.Range("M68").Value = 725
$ws.Range("N68").Value = "`$/unidad"
$ws.Range("O68").Value = "Provincia de Diguillín"
$ws.Range("P68").Value = 725
$ws.Range("Q68").Value = 1
$ws.Range("R68").Value = "Hortaliza"
